# Apply the "nuevos experimentos no convexos" numeric update across the
# non-convex bilevel experiment workbook.
#
# Helper: write a value into a cell, forcing it to be stored as TEXT
# (shared-string) even when the text looks like a number (e.g. "-2.9"),
# matching the original authoring tool's output. A leading apostrophe
# forces Excel to treat the value as text; we then reset the cell style
# back to Normal so no visible/persisted style difference remains.
function Set-TextValue {
    param(
        $Range,
        [string]$Text
    )
    $Range.Value = "'" + $Text
    $Range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# NOTE: worksheet names "Vector_bf" and "Vector_BF" differ only by case, and
# Worksheets.Item(<name>) lookup here is case-insensitive, so we address all
# sheets positionally (matches workbook.xml sheet order) to avoid ambiguity:
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha
$wsLider     = $wb.Worksheets.Item(2)
$wsFollower  = $wb.Worksheets.Item(3)
$wsPunto     = $wb.Worksheets.Item(4)
$wsVecbf     = $wb.Worksheets.Item(5)
$wsVecBF     = $wb.Worksheets.Item(6)

# --- Restricciones_del_lider (sheet: MIU_value constraints on x) ---
$wsLider.Range("A2").Value = "1.9 - x"
Set-TextValue $wsLider.Range("B2") "-2.9"
Set-TextValue $wsLider.Range("D2") "0.83"

$wsLider.Range("A3").Value = "-1.9 + x"
Set-TextValue $wsLider.Range("B3") "0.8999999999999999"
Set-TextValue $wsLider.Range("D3") "0.08"

# --- Restricciones_del_follower (sheet: Lambda/Beta/Gamma constraints on y) ---
$wsFollower.Range("A2").Value = "0.1499999999999999 - y"
Set-TextValue $wsFollower.Range("B2") "-1.15"
Set-TextValue $wsFollower.Range("D2") "0.01"
Set-TextValue $wsFollower.Range("E2") "0"
Set-TextValue $wsFollower.Range("F2") "0"

$wsFollower.Range("A3").Value = "-0.15000000000000002 + y"
Set-TextValue $wsFollower.Range("B3") "-0.85"
Set-TextValue $wsFollower.Range("D3") "0.97"
Set-TextValue $wsFollower.Range("E3") "0"
Set-TextValue $wsFollower.Range("F3") "0"

# --- Punto_modificado (x, y) ---
Set-TextValue $wsPunto.Range("A2") "1.9"
Set-TextValue $wsPunto.Range("B2") "0.15"

# --- Vector_bf ---
Set-TextValue $wsVecbf.Range("A2") "-0.39675000000000005"

# --- Vector_BF ---
Set-TextValue $wsVecBF.Range("A2") "1.6"
Set-TextValue $wsVecBF.Range("A3") "-4.9"
